$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert new sheet "PreInspection_Type" right after "MotorClaim_Insurer".
#    Copy MotorClaim_Insurer as a template: it already has the single
#    46.109375-wide column A plus the s=5/s=4/s=1 style pattern we need.
# ---------------------------------------------------------------------------
$motorClaim = $wb.Worksheets.Item("MotorClaim_Insurer")
$motorClaim.Copy($null, $motorClaim)
$preInsp = $wb.Worksheets.Item("MotorClaim_Insurer (2)")
$preInsp.Name = "PreInspection_Type"

# Trim to 7 rows / 1 column and drop the leftover sort state from the template.
$preInsp.Range("B1:B19").Clear()
$preInsp.Rows("8:19").Delete()
$preInsp.Sort.SortFields.Clear()

$preInsp.Range("A1").Value = "PRE INSPECTION TYPE (0,0)"
$preInsp.Range("A2").Value = "Car"
$preInsp.Range("A3").Value = "Car"
$preInsp.Range("A4").Value = "Commercial vehicles"
$preInsp.Range("A5").Value = "Construction equipment's"
$preInsp.Range("A6").Value = "Miscellaneous equipment's"
$preInsp.Range("A7").Value = "2-wheeler"

$preInsp.Range("A2").Validation.Modify(3, 1, 1, "=`$A`$3:`$A`$7")
$preInsp.Range("A3").Validation.Modify(3, 1, 1, "=`$A`$3:`$A`$7")

$preInsp.PageSetup.Orientation = 1
$preInsp.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2. Insert new sheet "Surveyor" right after "Repairer".
#    Copy Repairer as a template: matching A1:B4 layout, column widths and
#    s=5/s=3/s=4/s=1 style pattern.
# ---------------------------------------------------------------------------
$repairer = $wb.Worksheets.Item("Repairer")
$repairer.Copy($null, $repairer)
$surveyor = $wb.Worksheets.Item("Repairer (2)")
$surveyor.Name = "Surveyor"

$surveyor.Range("A1").Value = "SURVEYOR USERNAME (0,0)"
$surveyor.Range("A2").Value = "surveyorphase1@sendnow.win"
$surveyor.Range("A3").Value = "surveyorphase1@sendnow.win"
$surveyor.Range("A4").Value = "ehtasham@surveyororg.com"

$surveyor.Range("A2").Validation.Modify(3, 1, 1, "=`$A`$3:`$A`$4")
$surveyor.Range("B2").Validation.Modify(3, 1, 1, "=`$B`$3:`$B`$4")

$surveyor.Range("B8").Select()

# ---------------------------------------------------------------------------
# 3. MotorClaim_Insurer: selection becomes a full multi-cell range A1:A19.
# ---------------------------------------------------------------------------
$motorClaim.Range("A1:A19").Select()

# ---------------------------------------------------------------------------
# 4. SuperAdmin: no longer the active tab; selection moves to A11.
# ---------------------------------------------------------------------------
$superAdmin = $wb.Worksheets.Item("SuperAdmin")
$superAdmin.Range("A11").Select()

# ---------------------------------------------------------------------------
# 5. Customer: becomes the active tab; selection moves to K5; A2 value changes.
# ---------------------------------------------------------------------------
$customer = $wb.Worksheets.Item("Customer")
$customer.Range("A2").Value = "zooowlwhy@quicksend.ch"
$customer.Activate()
$customer.Range("K5").Select()

# ---------------------------------------------------------------------------
# 6. Repairer: selection moves to B18.
# ---------------------------------------------------------------------------
$repairer.Range("B18").Select()

# ---------------------------------------------------------------------------
# 7. Insurer_Log: selection moves to B18.
# ---------------------------------------------------------------------------
$insurerLog = $wb.Worksheets.Item("Insurer_Log")
$insurerLog.Range("B18").Select()

# Re-activate Customer last so it is the workbook's active/visible tab.
$customer.Activate()
